$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47; this shifts the existing rows 47-113
# down to 48-114 (preserving all their data/formatting), and leaves a
# blank row 47 ready to be populated with the new weekly record.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data record.
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44721
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100108
$ws.Range("H47").Value = "Tropicales y subtropicales"
$ws.Range("I47").Value = 100108002
$ws.Range("J47").Value = "Mango"
$ws.Range("K47").Value = "Sin especificar"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 180
$ws.Range("N47").Value = 9000
$ws.Range("O47").Value = 10000
$ws.Range("P47").Value = 9444
$ws.Range("Q47").Value = "$/bandeja 4 kilos"
$ws.Range("R47").Value = "Ecuador"
$ws.Range("S47").Value = 2361
$ws.Range("T47").Value = 4
